$d = $word.ActiveDocument

# 1) Remove the obsolete bullet items (duplicated "Minimum grain size..."
#    and "Check twin thickness..." entries plus "Randomly shift twin
#    placement" and "Make sure twins are separated enough"). Delete from
#    the bottom up so earlier paragraph indices stay valid.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -eq "Minimum grain size to insert twin user input (or rule)`r" -or
        $t -eq "Check twin thickness per grain to make sure its thicker than 1 voxel`r" -or
        $t -eq "Randomly shift twin placement`r" -or
        $t -eq "Make sure twins are separated enough`r") {
        $p.Range.Delete()
    }
}

# 2) Update the text of the remaining four bullet items.
[void]$d.Content.Find.Execute("Number of twins per grain user input", $true, $false, $false, $false, $false, $true, 1, $false, "Minimum grain size to insert twin user input (or rule)", 2)
[void]$d.Content.Find.Execute("Coherency fraction user input", $true, $false, $false, $false, $false, $true, 1, $false, "Check twin thickness per grain to make sure its thicker than 1 voxel", 2)
[void]$d.Content.Find.Execute("Twin spacing user input", $true, $false, $false, $false, $false, $true, 1, $false, "Add ensemble phase type", 2)
[void]$d.Content.Find.Execute("Twin morphology fraction user input (i.e. isthmus vs peninsula)", $true, $false, $false, $false, $false, $true, 1, $false, "Vary twin thickness", 2)

# 3) Move the "_GoBack" bookmark so that it now sits at the start of the
#    "Add ensemble phase type" paragraph (right before its run), matching
#    where the edit was last made.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Add ensemble phase type`r") {
        $r = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $r)
        break
    }
}
